$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(88, 3).Value = "Find_Drawdowns_test1"
$ws.Cells.Item(88, 2).Value = "Test find drawdown for simple returns"
$ws.Cells.Item(88, 1).Value = "Find_Drawdowns1"

$ws.Cells.Item(89, 1).Value = "Find_Drawdowns2"
$ws.Cells.Item(89, 2).Value = "Test find drawdown for compound returns"
$ws.Cells.Item(89, 3).Value = "Find_Drawdowns_test2"

$ws.Range("C89").Select()
